# Auto-generated edit script: updates crypto price/volume table to match latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cell, far outside the used range, used only as a "General"-format donor so that
# forcing a numeric-looking price string to stay text (via NumberFormat "@") does not leave
# a new/changed style on the edited cell (PasteSpecial formats resets it back to the default style).
$formatDonor = $ws.Range("ZZ1000")

$ws.Range("D2").Value = "34.473.51"
$ws.Range("E2").Value = "  -0.39%  "
$ws.Range("D3").Value = "1.804.55"
$ws.Range("E3").Value = "  -0.95%  "
$ws.Range("E4").Value = "  +0.28%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.78"
$formatDonor.Copy()
$ws.Range("D5").PasteSpecial(-4122)
$ws.Range("E5").Value = "  -0.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.576"
$formatDonor.Copy()
$ws.Range("D6").PasteSpecial(-4122)
$ws.Range("E6").Value = "  +2.91%  "
$ws.Range("E7").Value = "  +0.22%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "36.31"
$formatDonor.Copy()
$ws.Range("D8").PasteSpecial(-4122)
$ws.Range("E8").Value = "  +3.76%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.298"
$formatDonor.Copy()
$ws.Range("D9").PasteSpecial(-4122)
$ws.Range("E9").Value = "  -0.52%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0692"
$formatDonor.Copy()
$ws.Range("D10").PasteSpecial(-4122)
$ws.Range("E10").Value = "  -0.47%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0963"
$formatDonor.Copy()
$ws.Range("D11").PasteSpecial(-4122)
$ws.Range("E11").Value = "  +1.28%  "
$ws.Range("D12").Value = "2.063.78"
$ws.Range("E12").Value = "  -0.97%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.58"
$formatDonor.Copy()
$ws.Range("D13").PasteSpecial(-4122)
$ws.Range("E13").Value = "  +1.32%  "
$ws.Range("D14").Value = "1.805.02"
$ws.Range("E14").Value = "  -0.93%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.645"
$formatDonor.Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("E15").Value = "  -0.11%  "
$ws.Range("E16").Value = "  +2.61%  "
$ws.Range("D17").Value = "34.427.94"
$ws.Range("E17").Value = "  -0.42%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "70.16"
$formatDonor.Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("E18").Value = "  +1.28%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "245.25"
$formatDonor.Copy()
$ws.Range("D19").PasteSpecial(-4122)
$ws.Range("E19").Value = "  -1.06%  "
$ws.Range("D20").Value = "0.0₃0790"
$ws.Range("E20").Value = "  -1.51%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.58"
$formatDonor.Copy()
$ws.Range("D21").PasteSpecial(-4122)
$ws.Range("E21").Value = "  +0.25%  "
$ws.Range("E22").Value = "  +0.26%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.19"
$formatDonor.Copy()
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.16"
$formatDonor.Copy()
$ws.Range("D24").PasteSpecial(-4122)
$ws.Range("E24").Value = "  +3.89%  "
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "173.03"
$formatDonor.Copy()
$ws.Range("D25").PasteSpecial(-4122)
$ws.Range("E25").Value = "  +1.16%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.99"
$formatDonor.Copy()
$ws.Range("D26").PasteSpecial(-4122)
$ws.Range("E26").Value = "  +8.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.97"
$formatDonor.Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").Value = "  +0.74%  "
$ws.Range("E28").Value = "  +1.23%  "
$ws.Range("E29").Value = "  +0.15%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.02"
$formatDonor.Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("E30").Value = "  -0.32%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.85"
$formatDonor.Copy()
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("E31").Value = "  -0.42%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0529"
$formatDonor.Copy()
$ws.Range("D32").PasteSpecial(-4122)
$ws.Range("E32").Value = "  -0.64%  "
$ws.Range("E33").Value = "  -0.38%  "
$ws.Range("E34").Value = "  -2.42%  "
$ws.Range("D35").Value = "1.391.78"
$ws.Range("E35").Value = "  -1.75%  "
$ws.Range("E36").Value = "  -1.11%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.45"
$formatDonor.Copy()
$ws.Range("D37").PasteSpecial(-4122)
$ws.Range("E37").Value = "  -6.57%  "
$ws.Range("E38").Value = "  -0.55%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0190"
$formatDonor.Copy()
$ws.Range("D39").PasteSpecial(-4122)
$ws.Range("E39").Value = "  -0.93%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "82.71"
$formatDonor.Copy()
$ws.Range("D40").PasteSpecial(-4122)
$ws.Range("E40").Value = "  -4.45%  "
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.960"
$formatDonor.Copy()
$ws.Range("D41").PasteSpecial(-4122)
$ws.Range("E41").Value = "  +0.31%  "
$ws.Range("B42").Value = "MXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.83"
$formatDonor.Copy()
$ws.Range("D42").PasteSpecial(-4122)
$ws.Range("E42").Value = "  -0.81%  "
$ws.Range("E43").Value = "  +0.69%  "
$ws.Range("E44").Value = "  +8.25%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.45"
$formatDonor.Copy()
$ws.Range("D45").PasteSpecial(-4122)
$ws.Range("E45").Value = "  -3.78%  "
$ws.Range("E46").Value = "  -1.10%  "
$ws.Range("E47").Value = "  -4.22%  "
$ws.Range("D48").Value = "1.965.47"
$ws.Range("E48").Value = "  -1.08%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "104.37"
$formatDonor.Copy()
$ws.Range("D49").PasteSpecial(-4122)
$ws.Range("E49").Value = "  -1.64%  "
$ws.Range("E50").Value = "  +0.23%  "
$ws.Range("E51").Value = "  -1.97%  "

$excel.CutCopyMode = $false

